# Edit script for Parameters/ShortDivisionNames.xlsx
# Implements the changes described by the commit:
#   Added new settings values:
#    - employee quarterly rates;
#    - names of functional groups;
#    - hourly production calendar;
#    - other.
#
# Concretely (after resolving sharedStrings re-indexing noise in the raw
# OOXML diff) the only real content changes on sheet "Лист2" are:
#   B6 : "НПИС" -> "ДПИС"
#   A16:D16 : new row  ООО "АрСи БиАй" | RCBI            | 0 | -5
#   A17:D17 : new row  Департамент поддержки информационных систем | ДПИС | 0 | -1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист2")

# 1. Fill in the two previously-empty rows with new data, copying the
#    formatting used by the row directly above them (row 15) so the
#    new cells pick up the same styles (data style for A/B/D, percent
#    style for C) instead of the empty-row style.
$ws.Range("A15:D15").Copy() | Out-Null
$ws.Range("A16:D16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A17:D17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 'ООО "АрСи БиАй"'
$ws.Range("B16").Value = "RCBI"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = -5

$ws.Range("B17").Value = "ДПИС"
$ws.Range("A17").Value = "Департамент поддержки информационных систем"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = -1

# 2. Rename the short name of the existing "support division" row to
#    reuse the same "ДПИС" text introduced above.
$ws.Range("B6").Value = "ДПИС"

# 3. Cosmetic: move the active selection to A7, matching the saved view
#    state recorded in the workbook after the edit.
$ws.Range("A7").Select() | Out-Null
